$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G6").Value = 20240110
$ws.Range("B7").Value = "CQ"
$ws.Range("C7").Select() | Out-Null
